$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace one paragraph's full content (including any proofErr
# wrapped run(s)) with freshly authored WordprocessingML that no longer
# carries the now-stale grammar-check annotations. InsertXML must be scoped
# to the exact Range whose contents should change, so we target the
# paragraph's own Range (pilcrow included) and supply a full pkg:package
# WordOpenXML document containing a single replacement <w:p>.
# ---------------------------------------------------------------------------
function Set-ParagraphXml($paragraph, [string]$innerParagraphXml) {
    $range = $paragraph.Range
    $prefix = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>'
    $suffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $xml = $prefix + $innerParagraphXml + $suffix
    $range.InsertXML($xml)
}

# 1) "Que es un combo?" -> "¿Qué es un combo?" (fixed grammar, drop proofErr wrap)
$p3 = $d.Paragraphs(3)
$p3xml = '<w:p w:rsidR="00E56AE6" w:rsidRDefault="00E56AE6" w:rsidP="00E56AE6">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr><w:t>¿Qué es un combo?</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p3 $p3xml

# 2) "Porque lo quiero en mi juego?" -> "¿Por qué lo quiero en mi juego?" (drop proofErr wrap)
$p4 = $d.Paragraphs(4)
$p4xml = '<w:p w:rsidR="00E56AE6" w:rsidRDefault="00E56AE6" w:rsidP="00E56AE6">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr><w:t>¿Por qué lo quiero en mi juego?</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p4 $p4xml

# 3) Split "Combos difíciles y combos fáciles" into two runs so the (moved)
#    _GoBack bookmark can sit between them, keeping the trailing
#    " y como afectan la accesibilidad" run untouched.
$p7 = $d.Paragraphs(7)
$p7xml = '<w:p w:rsidR="00E56AE6" w:rsidRDefault="00E56AE6" w:rsidP="00E56AE6">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve">Combos difíciles y </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr><w:t>combos fáciles</w:t></w:r>' +
    '<w:r w:rsidR="001D54FB"><w:rPr><w:sz w:val="28"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve"> y como afectan la accesibilidad</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p7 $p7xml

# 4) Move the "_GoBack" bookmark from the end of the document to sit
#    between "Combos difíciles y " and "combos fáciles" - Word keeps a
#    single "_GoBack" bookmark, so re-adding it here removes the old one
#    automatically (mirrors Word's "last edit location" bookmark behaviour).
$p7 = $d.Paragraphs(7)
$searchRange = $p7.Range.Duplicate
$found = $searchRange.Find.Execute("combos fáciles")
$bookmarkRange = $d.Range($searchRange.Start, $searchRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
